# The document ends with two empty paragraphs. The very last one gets
# replaced (in place) by the whole new "resignation letter to Vence" -
# InsertXML on that paragraph's Range swaps its contents for the full run
# of new <w:p> siblings supplied below, in one shot.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$apos = [char]0x2019

$newParagraphs = @(
    "<w:p $wNs><w:r><w:t>Dear Vence,</w:t></w:r></w:p>"
    "<w:p $wNs><w:r><w:t>I appreciate the job opportunity you gave me.</w:t></w:r></w:p>"
    "<w:p $wNs><w:r><w:t>I learned a lot from it, and it was quite unlike any other job I ever had.</w:t></w:r></w:p>"
    "<w:p $wNs><w:r><w:t>However, I no longer feel comfortable gambling on the life of my co-workers every day.</w:t></w:r></w:p>"
    "<w:p $wNs><w:r><w:t>Also, I don${apos}t enjoy working for free.</w:t></w:r></w:p>"
    "<w:p $wNs><w:r><w:t>Therefore, please accept this letter as my official resignation, effective today.</w:t></w:r></w:p>"
    "<w:p $wNs><w:r><w:lastRenderedPageBreak/><w:t>May your bad days be the best days,</w:t></w:r></w:p>"
    "<w:p $wNs><w:r><w:t>You</w:t></w:r></w:p>"
) -join ""

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertXML($newParagraphs) | Out-Null
